$wb = $excel.ActiveWorkbook

# The workbook's 日期 (date) column used a short ROC-calendar numeric string
# (e.g. "11426" = year 114, month 2, day 6). The edit reformats every one of
# those date codes to a fixed-width YYYMMDD form (e.g. "1140206").
# These values are plain text (not real dates/numbers), so we force the
# cell to text format before writing the digit string, then restore the
# cell's style to Normal so only the *value* changes, not the formatting.
function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Map of sheet name -> ordered list of (cellRef, newValue) for column A.
$sheetEdits = @{
    "陽孟青2月報表" = @(
        @("A2", "1140206"),
        @("A3", "1140208"),
        @("A4", "1140224"),
        @("A5", "1140228"),
        @("A6", "1140301"),
        @("A7", "1140302"),
        @("A8", "1140303"),
        @("A9", "1140303"),
        @("A10", "1140304"),
        @("A11", "1140305"),
        @("A12", "1140306"),
        @("A13", "1140307"),
        @("A14", "1140308"),
        @("A15", "1140310"),
        @("A16", "1140311"),
        @("A17", "1140311")
    )
    "林李勝2月報表" = @(
        @("A2", "1140206"),
        @("A3", "1140208"),
        @("A4", "1140224"),
        @("A5", "1140228"),
        @("A6", "1140301"),
        @("A7", "1140301"),
        @("A8", "1140303"),
        @("A9", "1140304"),
        @("A10", "1140305"),
        @("A11", "1140306"),
        @("A12", "1140307"),
        @("A13", "1140308"),
        @("A14", "1140310"),
        @("A15", "1140311")
    )
    "王登翊2月報表" = @(
        @("A2", "1140205"),
        @("A3", "1140206"),
        @("A4", "1140208"),
        @("A5", "1140218"),
        @("A6", "1140224"),
        @("A7", "1140226"),
        @("A8", "1140228"),
        @("A9", "1140301"),
        @("A10", "1140302"),
        @("A11", "1140303"),
        @("A12", "1140303"),
        @("A13", "1140304"),
        @("A14", "1140304"),
        @("A15", "1140305"),
        @("A16", "1140306"),
        @("A17", "1140307"),
        @("A18", "1140307"),
        @("A19", "1140308"),
        @("A20", "1140310"),
        @("A21", "1140311")
    )
    "張浩宇2月報表" = @(
        @("A2", "1140206"),
        @("A3", "1140207"),
        @("A4", "1140208"),
        @("A5", "1140218"),
        @("A6", "1140224"),
        @("A7", "1140228"),
        @("A8", "1140303"),
        @("A9", "1140303"),
        @("A10", "1140304"),
        @("A11", "1140305"),
        @("A12", "1140306"),
        @("A13", "1140306"),
        @("A14", "1140307"),
        @("A15", "1140307"),
        @("A16", "1140311"),
        @("A17", "1140311"),
        @("A18", "1140312")
    )
    "詹智勝3月報表" = @(
        @("A2", "1140301"),
        @("A3", "1140303"),
        @("A4", "1140304"),
        @("A5", "1140306"),
        @("A6", "1140307"),
        @("A7", "1140307"),
        @("A8", "1140308"),
        @("A9", "1140311"),
        @("A10", "1140311"),
        @("A11", "1140311"),
        @("A12", "1140312")
    )
    "林長億2月報表" = @(
        @("A2", "1140208"),
        @("A3", "1140218"),
        @("A4", "1140224"),
        @("A5", "1140228"),
        @("A6", "1140301"),
        @("A7", "1140301"),
        @("A8", "1140302"),
        @("A9", "1140303"),
        @("A10", "1140312")
    )
    "張國祥2月報表" = @(
        @("A2", "1140206"),
        @("A3", "1140207"),
        @("A4", "1140208"),
        @("A5", "1140218"),
        @("A6", "1140305"),
        @("A7", "1140306"),
        @("A8", "1140307"),
        @("A9", "1140308"),
        @("A10", "1140309"),
        @("A11", "1140312"),
        @("A12", "1140315")
    )
}

foreach ($sheetName in $sheetEdits.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($edit in $sheetEdits[$sheetName]) {
        $cellRef = $edit[0]
        $newValue = $edit[1]
        Set-TextValue $ws.Range($cellRef) $newValue
    }
}
